$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.075936913490295
$ws.Range("B1").Value = 2.319365978240967
$ws.Range("C1").Value = 9.471441268920898
$ws.Range("D1").Value = 2.283515214920044
$ws.Range("E1").Value = 1.309595704078674
